$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert 7 new rows after the existing last data row (19) to make room
# for the expanded worker table (rows 16-26) while preserving the signature
# block spacing below it.
$ws.Rows("20:26").Insert()

# Step 2: Re-apply the "middle" row formatting (no bottom border) to rows 16-25,
# and the "last row" formatting (bottom border) to row 26 (now the final worker row).
# NOTE: copy the "last row" style (from the original row 19) onto row 26 BEFORE
# overwriting row 19 itself with the "middle" style, otherwise row 19's special
# formatting is lost before it can be copied.
$ws.Range("B19:J19").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: Write the refreshed worker rows (11 total).
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "1143352669"
$ws.Cells.Item(16, 4).Value = "CARLOS MARIO FUENTES MAYA"
$ws.Cells.Item(16, 5).Value = "2211"
$ws.Cells.Item(16, 6).Value = 105620
$ws.Cells.Item(16, 7).Value = 2640500

$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "1047451044"
$ws.Cells.Item(17, 4).Value = "EVELIN ZAPATEIRO MARQUEZ"
$ws.Cells.Item(17, 5).Value = "2505"
$ws.Cells.Item(17, 6).Value = 4615
$ws.Cells.Item(17, 7).Value = 3461511

$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "1047451044"
$ws.Cells.Item(18, 4).Value = "EVELIN ZAPATEIRO MARQUEZ"
$ws.Cells.Item(18, 5).Value = "2503"
$ws.Cells.Item(18, 6).Value = 4234
$ws.Cells.Item(18, 7).Value = 3461511

$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "1047475389"
$ws.Cells.Item(19, 4).Value = "CARLOS FERNANDO OZUNA CORTINA"
$ws.Cells.Item(19, 5).Value = "2008"
$ws.Cells.Item(19, 6).Value = 40000
$ws.Cells.Item(19, 7).Value = 1984811

$ws.Cells.Item(20, 2).Value = "CC"
$ws.Cells.Item(20, 3).Value = "1050952446"
$ws.Cells.Item(20, 4).Value = "MERY ELLEN ESPINOSA RHENAL"
$ws.Cells.Item(20, 5).Value = "2505"
$ws.Cells.Item(20, 6).Value = 8015
$ws.Cells.Item(20, 7).Value = 6011280

$ws.Cells.Item(21, 2).Value = "CC"
$ws.Cells.Item(21, 3).Value = "1235044752"
$ws.Cells.Item(21, 4).Value = "JOSE CARLOS RIOS MARQUEZ"
$ws.Cells.Item(21, 5).Value = "2008"
$ws.Cells.Item(21, 6).Value = 26666
$ws.Cells.Item(21, 7).Value = 5224500

$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "1143401581"
$ws.Cells.Item(22, 4).Value = "DONNA MARGARITA CANO RAMIREZ"
$ws.Cells.Item(22, 5).Value = "2507"
$ws.Cells.Item(22, 6).Value = 6067
$ws.Cells.Item(22, 7).Value = 4550175

$ws.Cells.Item(23, 2).Value = "CC"
$ws.Cells.Item(23, 3).Value = "1047496658"
$ws.Cells.Item(23, 4).Value = "ANDERSON JOSE PEREZ BELTRAN"
$ws.Cells.Item(23, 5).Value = "2108"
$ws.Cells.Item(23, 6).Value = 60000
$ws.Cells.Item(23, 7).Value = 1500000

$ws.Cells.Item(24, 2).Value = "CC"
$ws.Cells.Item(24, 3).Value = "1047496658"
$ws.Cells.Item(24, 4).Value = "ANDERSON JOSE PEREZ BELTRAN"
$ws.Cells.Item(24, 5).Value = "2107"
$ws.Cells.Item(24, 6).Value = 60000
$ws.Cells.Item(24, 7).Value = 1500000

$ws.Cells.Item(25, 2).Value = "CC"
$ws.Cells.Item(25, 3).Value = "1047496658"
$ws.Cells.Item(25, 4).Value = "ANDERSON JOSE PEREZ BELTRAN"
$ws.Cells.Item(25, 5).Value = "2106"
$ws.Cells.Item(25, 6).Value = 60000
$ws.Cells.Item(25, 7).Value = 1500000

$ws.Cells.Item(26, 2).Value = "CE"
$ws.Cells.Item(26, 3).Value = "20394544"
$ws.Cells.Item(26, 4).Value = "LUZ ELENA UTRIA ORTIZ"
$ws.Cells.Item(26, 5).Value = "2210"
$ws.Cells.Item(26, 6).Value = 9333
$ws.Cells.Item(26, 7).Value = 1000000

# Step 4: Update the summary figures at the top of the statement.
$ws.Range("E11").Value = 384550
$ws.Range("C13").Value = 8
$ws.Range("F13").Value = 9
